# Working hours workbook: add three new logged shifts (2014-04-05 x2, 2014-04-07)
# just above the trailing blank/summary rows, pushing the summary block down
# and refreshing its totals.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet currently ends with:
#   row 116: blank separator row (D/E/F styled but empty)
#   row 117: "sum [min]"
#   row 118: "sum [h]"
#   row 119: "sum [working weeks]"
# Insert 3 fresh rows right before that blank separator row so the
# separator + summary rows shift down to 119-122, and the new data rows
# become 116-118.
$ws.Rows("116:118").Insert() | Out-Null

# --- Row 116: 2014-04-05, 09:30 -> 12:45 ---
$ws.Range("A116").Value2 = 2014
$ws.Range("B116").Value2 = 4
$ws.Range("C116").Value2 = 5
$ws.Range("D116").Value2 = 0.39583333333333331
$ws.Range("E116").Value2 = 0.53125

# --- Row 117: 2014-04-05, 13:00 -> 15:45 ---
$ws.Range("A117").Value2 = 2014
$ws.Range("B117").Value2 = 4
$ws.Range("C117").Value2 = 5
$ws.Range("D117").Value2 = 0.54166666666666663
$ws.Range("E117").Value2 = 0.65625

# --- Row 118: 2014-04-07, 19:30 -> 20:00 ---
$ws.Range("A118").Value2 = 2014
$ws.Range("B118").Value2 = 4
$ws.Range("C118").Value2 = 7
$ws.Range("D118").Value2 = 0.8125
$ws.Range("E118").Value2 = 0.83333333333333337

# Carry the "time spent [min]" / "time spent [h]" formulas down onto the
# three new rows, matching the pattern used for every other data row.
$ws.Range("F116:F118").Formula = "=(E116-D116)*24*60"
$ws.Range("G116:G118").Formula = "=F116/60"

# Keep the view roughly where the author left it (near the bottom of the
# now-longer sheet), with J113 as the active cell.
$excel.ActiveWindow.ScrollRow = 85
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("J113").Select() | Out-Null
